# Modification noms 2 variables + relecture
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header labels (column F and column J, row 1)
$ws.Range("F1").Value = "Part non importée volailles (%)"
$ws.Range("J1").Value = "Part non perdue (%)"

# Column F now represents the complementary share (100 - old import %)
# for each data row (rows 2 through 29).
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 6)  # column F
    $oldValue = $cell.Value2
    $cell.Value = 100 - $oldValue
}
